# ERS requisitos de agente de playa en raiz
# Rename sheets, reorganize active view, and populate the "tipos infraccion"
# sheet with the list of unique "Tipo" values (and their source row positions)
# taken from the "codigos" sheet.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)

# 1) Rename the first two sheets.
$ws1.Name = "codigos"
$ws2.Name = "tipos infraccion"

# 2) Populate "tipos infraccion" (ws2) column A with the distinct "Tipo"
#    values pulled from "codigos" (ws1) column A, placed at the rows where
#    each value first begins in the source sheet layout.
$srcRows = @(1,2,5,8,14,16,21,25,27,34,35,45,47,49,52,53,54,75,77,99,114,126,130,134)
$dstRows = @(1,2,3,4,5,6,7,8,9,10,11,12,13,14,15,16,17,38,40,62,77,89,93,97)

for ($i = 0; $i -lt $srcRows.Length; $i++) {
    $srcRow = $srcRows[$i]
    $dstRow = $dstRows[$i]
    $ws2.Cells.Item($dstRow, 1).Value = $ws1.Cells.Item($srcRow, 1).Value()
}

# 3) Copy formatting: row 1 is the header style, the rest share the normal
#    "Tipo" column style already used throughout "codigos".
$ws1.Range("A1").Copy()
$ws2.Range("A1").PasteSpecial(-4122) | Out-Null

$ws1.Range("A2").Copy()
$ws2.Range("A2:A17").PasteSpecial(-4122) | Out-Null
$ws2.Range("A38").PasteSpecial(-4122) | Out-Null
$ws2.Range("A40").PasteSpecial(-4122) | Out-Null
$ws2.Range("A62").PasteSpecial(-4122) | Out-Null
$ws2.Range("A77").PasteSpecial(-4122) | Out-Null
$ws2.Range("A89").PasteSpecial(-4122) | Out-Null
$ws2.Range("A93").PasteSpecial(-4122) | Out-Null
$ws2.Range("A97").PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Match column width used on "codigos" column A.
$ws2.Range("A1").ColumnWidth = 27.7109375

# 4) Update selections / active sheet so that "tipos infraccion" becomes the
#    active tab, with the view scrolled similarly to the source workbook.
$ws1.Range("C47").Select() | Out-Null
$ws2.Activate() | Out-Null
$ws2.Range("A16").Select() | Out-Null
